# Weekly refresh of Fruta/Hortaliza figures for "Agricola del Norte S.A. de
# Arica - Haba" (Hortaliza). The date (D), Volumen (J), Precio minimo (K),
# Precio maximo (L), Precio promedio ponderado (M) and Precio $/Kg (P)
# columns are re-shuffled row by row for this week's refresh; every other
# column is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => Fecha(serial), Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg
$rows = @(
    @{ Row = 2;  D = 44607; J = 900;  K = 1300; L = 1400; M = 1350; P = 1350 },
    @{ Row = 3;  D = 44284; J = 1500; K = 800;  L = 850;  M = 825;  P = 825  },
    @{ Row = 4;  D = 44656; J = 1000; K = 900;  L = 1000; M = 950;  P = 950  },
    @{ Row = 5;  D = 44229; J = 1500; K = 1400; L = 1500; M = 1450; P = 1450 },
    @{ Row = 6;  D = 44455; J = 1100; K = 900;  L = 1000; M = 950;  P = 950  },
    @{ Row = 7;  D = 44407; J = 1000; K = 1200; L = 1300; M = 1250; P = 1250 },
    @{ Row = 8;  D = 44649; J = 600;  K = 900;  L = 1000; M = 950;  P = 950  },
    @{ Row = 9;  D = 44291; J = 1000; K = 1000; L = 1200; M = 1100; P = 1100 },
    @{ Row = 10; D = 44341; J = 1300; K = 900;  L = 1000; M = 950;  P = 950  },
    @{ Row = 11; D = 44175; J = 1600; K = 1000; L = 1200; M = 1100; P = 1100 },
    @{ Row = 12; D = 44449; J = 1300; K = 900;  L = 950;  M = 925;  P = 925  },
    @{ Row = 13; D = 44476; J = 900;  K = 700;  L = 800;  M = 750;  P = 750  },
    @{ Row = 14; D = 44687; J = 1000; K = 1200; L = 1300; M = 1250; P = 1250 },
    @{ Row = 15; D = 44550; J = 1300; K = 1000; L = 1200; M = 1100; P = 1100 },
    @{ Row = 16; D = 44243; J = 1200; K = 1200; L = 1300; M = 1250; P = 1250 },
    @{ Row = 17; D = 44442; J = 1250; K = 850;  L = 900;  M = 875;  P = 875  },
    @{ Row = 18; D = 44638; J = 1000; K = 900;  L = 950;  M = 925;  P = 925  },
    @{ Row = 19; D = 44484; J = 900;  K = 750;  L = 800;  M = 775;  P = 775  },
    @{ Row = 20; D = 44784; J = 1000; K = 1200; L = 1300; M = 1250; P = 1250 },
    @{ Row = 21; D = 44673; J = 900;  K = 1300; L = 1400; M = 1350; P = 1350 },
    @{ Row = 22; D = 44453; J = 1000; K = 800;  L = 900;  M = 850;  P = 850  }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($r.Row, 10).Value = $r.J   # J: Volumen
    $ws.Cells.Item($r.Row, 11).Value = $r.K   # K: Precio minimo
    $ws.Cells.Item($r.Row, 12).Value = $r.L   # L: Precio maximo
    $ws.Cells.Item($r.Row, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # P: Precio $/Kg
}
